$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pre-Alert Template Import")

# Cells A3, B3, C3, AN3, AO3 use a quote-prefixed text style (forced text),
# so set via Formula with a leading apostrophe to preserve that style.
$ws.Range("A3").Formula = "'JSSO1000250"
$ws.Range("B3").Formula = "'JSSO1000250"
$ws.Range("C3").Formula = "'JSSO1000250"
$ws.Range("AN3").Formula = "'MBLJSSO1000250"
$ws.Range("AO3").Formula = "'HBLJSSO1000250"

# AJ3 uses a plain text style (no quote prefix needed).
$ws.Range("AJ3").Value = "JSCN1000250"
